$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily refresh of the 水桶/配送 tracking sheet:
#   column D = total cycle length (days)
#   column E = remaining days in the current cycle
#   column F = start date (yyyymmdd) of the current cycle
#
# Each day: E decrements by 1. Once a cycle's remaining days would drop
# to 0, the cycle restarts: E is reset to the full cycle length (D) and
# F is advanced by D days (the new cycle's start date).

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $total = $ws.Range("D$r").Value()
    $remaining = $ws.Range("E$r").Value()
    $startDate = $ws.Range("F$r").Value()

    if ($null -eq $total -or $null -eq $remaining -or $null -eq $startDate) {
        continue
    }

    # Every row's start date is parsed up front (it anchors the cycle
    # math below). A handful of rows carry a corrupted value that isn't
    # a proper 8-digit yyyymmdd (e.g. an extra stray digit) — those are
    # left completely untouched, same as the source automation.
    if ($startDate -lt 10000101 -or $startDate -gt 99991231) {
        continue
    }

    $year = [Math]::Floor($startDate / 10000)
    $month = [Math]::Floor(($startDate % 10000) / 100)
    $day = $startDate % 100

    $ok = $true
    try {
        $parsedDate = Get-Date -Year $year -Month $month -Day $day
    } catch {
        $ok = $false
    }

    if (-not $ok) {
        continue
    }

    $newRemaining = $remaining - 1

    if ($newRemaining -le 0) {
        # Cycle rollover: reset the remaining-days counter and advance
        # the start date by the full cycle length.
        $newDate = $parsedDate.AddDays($total)
        $newStartDate = [int]($newDate.ToString("yyyyMMdd"))

        $ws.Range("E$r").Value = $total
        $ws.Range("F$r").Value = $newStartDate
    } else {
        $ws.Range("E$r").Value = $newRemaining
    }
}
